$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "48.670.03"
$ws.Range("E2").Value = "  -2.06%  "

$ws.Range("D3").Value = "2.623.18"
$ws.Range("E3").Value = "  +2.93%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "109.67"
$ws.Range("E5").Value = "  +1.07%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "321.88"
$ws.Range("E6").Value = "  +0.05%  "

$ws.Range("E7").Value = "  -1.30%  "

$ws.Range("E8").Value = "  +0.11%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.537"
$ws.Range("E9").Value = "  -3.44%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.15"
$ws.Range("E10").Value = "  -2.45%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.78"
$ws.Range("E11").Value = "  -2.51%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0807"
$ws.Range("E12").Value = "  -0.78%  "

$ws.Range("E13").Value = "  +0.21%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.19"
$ws.Range("E14").Value = "  -0.67%  "

$ws.Range("D15").Value = "3.036.09"
$ws.Range("E15").Value = "  +3.13%  "

$ws.Range("D16").Value = "2.634.80"
$ws.Range("E16").Value = "  +1.86%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.859"
$ws.Range("E17").Value = "  +0.35%  "

$ws.Range("D18").Value = "48.626.41"
$ws.Range("E18").Value = "  -1.76%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.76"
$ws.Range("E19").Value = "  -3.27%  "

$ws.Range("B20").Value = "ImmutableX"
$ws.Range("C20").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.90"
$ws.Range("E20").Value = "  -4.27%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.66"
$ws.Range("E21").Value = "  -0.03%  "

$ws.Range("D22").Value = "0.0₃0939"
$ws.Range("E22").Value = "  -0.36%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "269.04"
$ws.Range("E23").Value = "  -6.25%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.63"
$ws.Range("E24").Value = "  -4.55%  "

$ws.Range("E25").Value = "  -0.19%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "25.97"
$ws.Range("E26").Value = "  -1.48%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("E27").Value = "  -0.01%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.03"
$ws.Range("E28").Value = "  +2.18%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.22"
$ws.Range("E29").Value = "  -0.21%  "

$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.137"
$ws.Range("E30").Value = "  -4.65%  "

$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.71"
$ws.Range("E31").Value = "  -1.42%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.29"
$ws.Range("E32").Value = "  -0.29%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.45"
$ws.Range("E33").Value = "  +1.87%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.22"
$ws.Range("E34").Value = "  -1.78%  "

$ws.Range("E35").Value = "  -0.09%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0792"
$ws.Range("E36").Value = "  +1.01%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.92"
$ws.Range("E37").Value = "  +5.80%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.02"
$ws.Range("E38").Value = "  +0.75%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.13"
$ws.Range("E39").Value = "  +4.86%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "125.33"
$ws.Range("E40").Value = "  +3.89%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.41"
$ws.Range("E41").Value = "  +0.26%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.110"
$ws.Range("E42").Value = "  -0.99%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.14"
$ws.Range("E43").Value = "  -3.18%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0313"
$ws.Range("E44").Value = "  +1.21%  "

$ws.Range("D45").Value = "2.064.36"
$ws.Range("E45").Value = "  +2.65%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.22"
$ws.Range("E46").Value = "  -1.73%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.11"
$ws.Range("E47").Value = "  +5.03%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.17"
$ws.Range("E48").Value = "  +1.43%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.93"
$ws.Range("E49").Value = "  -1.22%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "58.44"
$ws.Range("E50").Value = "  +2.12%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.14"
$ws.Range("E51").Value = "  -3.23%  "

